$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 105.384211
$ws.Range("H2").Value = 316.152633
$ws.Range("I2").Value = 0.1017338963855636
$ws.Range("J2").Value = 0.1017338963855636
$ws.Range("M2").Value = 8.142376
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 858.0778704253361
$ws.Range("R2").Value = 7722.700833828026
$ws.Range("S2").Value = 0.01771506512626762
$ws.Range("T2").Value = 0.01771506512626762
$ws.Range("G3").Value = 105.384211
$ws.Range("H3").Value = 316.152633
$ws.Range("I3").Value = 0.1017338963855636
$ws.Range("J3").Value = 0.1017338963855636
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 2565.087983036655
$ws.Range("R3").Value = 23085.79184732989
$ws.Range("S3").Value = 0.05295638337762577
$ws.Range("T3").Value = 0.05295638337762578
$ws.Range("G4").Value = 105.384211
$ws.Range("H4").Value = 316.152633
$ws.Range("I4").Value = 0.1017338963855636
$ws.Range("J4").Value = 0.1017338963855636
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 1504.595040352374
$ws.Range("R4").Value = 13541.35536317137
$ws.Range("S4").Value = 0.03106244788167022
$ws.Range("T4").Value = 0.03106244788167022
$ws.Range("I5").Value = 0.5920257690987943
$ws.Range("J5").Value = 0.5920257690987943
$ws.Range("M5").Value = 8.142376
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 4993.460677647875
$ws.Range("R5").Value = 44941.14609883088
$ws.Range("S5").Value = 0.1030902720590388
$ws.Range("T5").Value = 0.1030902720590388
$ws.Range("I6").Value = 0.5920257690987943
$ws.Range("J6").Value = 0.5920257690987943
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("S6").Value = 0.3081720519089289
$ws.Range("T6").Value = 0.3081720519089289
$ws.Range("I7").Value = 0.5920257690987943
$ws.Range("J7").Value = 0.5920257690987943
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("S7").Value = 0.1807634451308266
$ws.Range("T7").Value = 0.1807634451308266
$ws.Range("H8").Value = 951.685638
$ws.Range("I8").Value = 0.306240334515642
$ws.Range("J8").Value = 0.306240334515642
$ws.Range("M8").Value = 8.142376
$ws.Range("N8").Value = 24.427128
$ws.Range("O8").Value = 0.1741313933276368
$ws.Range("P8").Value = 0.1741313933276368
$ws.Range("Q8").Value = 2582.994099465296
$ws.Range("R8").Value = 23246.94689518767
$ws.Range("S8").Value = 0.05332605614233031
$ws.Range("T8").Value = 0.05332605614233032
$ws.Range("H9").Value = 951.685638
$ws.Range("I9").Value = 0.306240334515642
$ws.Range("J9").Value = 0.306240334515642
$ws.Range("O9").Value = 0.5205382400466131
$ws.Range("P9").Value = 0.5205382400466131
$ws.Range("R9").Value = 69493.06837802407
$ws.Range("S9").Value = 0.1594098047600583
$ws.Range("T9").Value = 0.1594098047600584
$ws.Range("H10").Value = 951.685638
$ws.Range("I10").Value = 0.306240334515642
$ws.Range("J10").Value = 0.306240334515642
$ws.Range("O10").Value = 0.3053303666257501
$ws.Range("P10").Value = 0.3053303666257501
$ws.Range("R10").Value = 40762.31564449588
$ws.Range("S10").Value = 0.09350447361325333
$ws.Range("T10").Value = 0.09350447361325334
